$d = $word.ActiveDocument

# Fix the incorrect script name in the doco: "mai_4300_fix_1.sql" -> "mai_4300_fix1.sql"
$old = "start mai_4300_fix_1.sql and then press return"
$new = "start mai_4300_fix1.sql and then press return"
$findRange = $d.Content
$findRange.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null

# Re-locate the corrected sentence so we can split the run where the cursor
# would naturally have been left after typing (this also relocates the
# "_GoBack" bookmark from the top of the document to that edit point, just
# like Word does automatically when you make an edit).
$hit = $d.Content
$hit.Find.Execute($new, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$splitPos = $hit.Start + 18
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
